$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 799.5
$ws.Range("I32").Value = 600
$ws.Range("J32").Value = 999
$ws.Range("K32").Value = 600
$ws.Range("L32").Value = 999
$ws.Range("M32").Value = -274
$ws.Range("N32").Value = -1651
$ws.Range("H38").Value = 2007
$ws.Range("I38").Value = 2007
$ws.Range("K38").Value = 6021
$ws.Range("M38").Value = -5649
$ws.Range("H40").Value = 2999.5
$ws.Range("J40").Value = 2999.3333
$ws.Range("L40").Value = 2999.3333
$ws.Range("N40").Value = -3349.3333
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("M42").ClearContents()
$ws.Range("N42").ClearContents()
$ws.Range("H51").Value = 5386
$ws.Range("J51").Value = 3998.8
$ws.Range("L51").Value = 3998.8
$ws.Range("N51").Value = -4966.8
$ws.Range("H53").Value = 7812.385
$ws.Range("I53").Value = 12518.125
$ws.Range("K53").Value = 12518.125
$ws.Range("M53").Value = -11881.125
$ws.Range("H64").Value = 4666.3335
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 4666.3335
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 4666.3335
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -5162.3335
$ws.Range("H67").Value = 4666.3335
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 4666.3335
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 4666.3335
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -6382.3335
$ws.Range("H100").Value = 2374.75
$ws.Range("I100").Value = 1750
$ws.Range("J100").Value = 2999.5
$ws.Range("K100").Value = 1750
$ws.Range("L100").Value = 2999.5
$ws.Range("M100").Value = -1209
$ws.Range("N100").Value = -4081.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2900.9028
$ws.Range("I32").Value = 1711.3771
$ws.Range("K32").Value = 1711.3771
$ws.Range("M32").Value = -1424.3771
$ws.Range("H92").Value = 45499.668
$ws.Range("J92").Value = 45499.668
$ws.Range("L92").Value = 45499.668
$ws.Range("N92").Value = -50491.668
$ws.Range("H109").Value = 67885.5
$ws.Range("J109").Value = 67885.5
$ws.Range("L109").Value = 67885.5
$ws.Range("N109").Value = -70659.5
$ws.Range("H135").Value = 39747.25
$ws.Range("J135").Value = 39747.25
$ws.Range("L135").Value = 39747.25
$ws.Range("N135").Value = -49887.25

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 6745
$ws.Range("I80").Value = 86.59999999999999
$ws.Range("J80").Value = 10074.2
$ws.Range("K80").Value = 86.59999999999999
$ws.Range("L80").Value = 10074.2
$ws.Range("M80").Value = 911.4
$ws.Range("N80").Value = -12070.2
$ws.Range("H83").Value = 6745
$ws.Range("I83").Value = 86.59999999999999
$ws.Range("J83").Value = 10074.2
$ws.Range("K83").Value = 433
$ws.Range("L83").Value = 50371
$ws.Range("M83").Value = 4559
$ws.Range("N83").Value = -60355

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 974.8333
$ws.Range("I16").Value = 929.8
$ws.Range("J16").Value = 1200
$ws.Range("K16").Value = 929.8
$ws.Range("L16").Value = 1200
$ws.Range("M16").Value = -642.8
$ws.Range("N16").Value = -1774
$ws.Range("H31").Value = 3132.72
$ws.Range("I31").Value = 1179.9333
$ws.Range("K31").Value = 1179.9333
$ws.Range("M31").Value = -884.9332999999999
$ws.Range("H34").Value = 3132.72
$ws.Range("I34").Value = 1179.9333
$ws.Range("K34").Value = 1179.9333
$ws.Range("M34").Value = -977.9332999999999
$ws.Range("H113").Value = 974.8333
$ws.Range("I113").Value = 929.8
$ws.Range("J113").Value = 1200
$ws.Range("K113").Value = 929.8
$ws.Range("L113").Value = 1200
$ws.Range("M113").Value = 1240.2
$ws.Range("N113").Value = -5540
$ws.Range("H132").Value = 2557.6428
$ws.Range("I132").Value = 1179.5
$ws.Range("K132").Value = 3538.5
$ws.Range("M132").Value = -1008.5
$ws.Range("H134").Value = 922.6667
$ws.Range("I134").Value = 887.4
$ws.Range("K134").Value = 2662.2
$ws.Range("M134").Value = -127.1999999999998

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 578.3570999999999
$ws.Range("I5").Value = 516
$ws.Range("K5").Value = 1548
$ws.Range("M5").Value = -1436
$ws.Range("H33").Value = 80.2
$ws.Range("J33").Value = 39.2
$ws.Range("L33").Value = 235.2
$ws.Range("N33").Value = -801.2
$ws.Range("H68").Value = 2000
$ws.Range("J68").Value = 2000
$ws.Range("L68").Value = 6000
$ws.Range("N68").Value = -7622
$ws.Range("H71").Value = 2000
$ws.Range("J71").Value = 2000
$ws.Range("L71").Value = 18000
$ws.Range("N71").Value = -26112
$ws.Range("H80").Value = 5109
$ws.Range("I80").Value = 5163.5
$ws.Range("K80").Value = 15490.5
$ws.Range("M80").Value = -14554.5
$ws.Range("H83").Value = 5109
$ws.Range("I83").Value = 5163.5
$ws.Range("K83").Value = 46471.5
$ws.Range("M83").Value = -41791.5
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").ClearContents()
$ws.Range("H107").Value = 575.7646999999999
$ws.Range("J107").Value = 575.7646999999999
$ws.Range("L107").Value = 1727.2941
$ws.Range("N107").Value = -5567.2941
$ws.Range("H122").Value = 897.4167
$ws.Range("J122").Value = 1579.25
$ws.Range("L122").Value = 14213.25
$ws.Range("N122").Value = -19113.25
$ws.Range("H129").Value = 39054.42
$ws.Range("J129").Value = 52799.145
$ws.Range("L129").Value = 158397.435
$ws.Range("N129").Value = -168397.435
$ws.Range("H132").Value = 1627
$ws.Range("J132").Value = 1761.6
$ws.Range("L132").Value = 15854.4
$ws.Range("N132").Value = -20914.4
$ws.Range("H135").Value = 578.3570999999999
$ws.Range("I135").Value = 516
$ws.Range("K135").Value = 4644
$ws.Range("M135").Value = -2109

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1399.7
$ws.Range("I113").Value = 1193.4
$ws.Range("J113").Value = 1606
$ws.Range("K113").Value = 1193.4
$ws.Range("L113").Value = 1606
$ws.Range("M113").Value = 976.5999999999999
$ws.Range("N113").Value = -5946
$ws.Range("H122").Value = 2048.55
$ws.Range("I122").Value = 1987.5454
$ws.Range("K122").Value = 5962.6362
$ws.Range("M122").Value = -3512.6362

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1574.6666
$ws.Range("I22").Value = 1533.5
$ws.Range("J22").Value = 1602.1111
$ws.Range("K22").Value = 1533.5
$ws.Range("L22").Value = 1602.1111
$ws.Range("M22").Value = -1238.5
$ws.Range("N22").Value = -2192.1111
$ws.Range("H27").Value = 1574.6666
$ws.Range("I27").Value = 1533.5
$ws.Range("J27").Value = 1602.1111
$ws.Range("K27").Value = 1533.5
$ws.Range("L27").Value = 1602.1111
$ws.Range("M27").Value = -1426.5
$ws.Range("N27").Value = -1816.1111
$ws.Range("H46").Value = 1891.5834
$ws.Range("J46").Value = 1945.3636
$ws.Range("L46").Value = 1945.3636
$ws.Range("N46").Value = -2321.3636
$ws.Range("H55").Value = 326.92856
$ws.Range("I55").Value = 104.375
$ws.Range("J55").Value = 623.6667
$ws.Range("K55").Value = 104.375
$ws.Range("L55").Value = 623.6667
$ws.Range("M55").Value = 68.625
$ws.Range("N55").Value = -969.6667
$ws.Range("H132").Value = 2195.7058
$ws.Range("I132").Value = 1535.9
$ws.Range("J132").Value = 3138.2856
$ws.Range("K132").Value = 4607.700000000001
$ws.Range("L132").Value = 9414.856800000001
$ws.Range("M132").Value = -2077.700000000001
$ws.Range("N132").Value = -14474.8568
